$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: fiscal period headers (shift + append new 1401 period)
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# Row 9: publish-date headers (shift + append new date)
$ws.Range("D9").Value = "1399-03-13 (8)"
$ws.Range("E9").Value = "1400-03-11 (10)"
$ws.Range("F9").Value = "1401-03-24 (10)"
$ws.Range("G9").Value = "1402-02-28 (8)"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "1402-02-28"

# Rows 11-27: shift financial data left one column, append new rightmost column
# Row 11
$ws.Range("D11").Value = 2090302
$ws.Range("E11").Value = 3554839
$ws.Range("F11").Value = 7290461
$ws.Range("G11").Value = 9566926
$ws.Range("H11").Value = 13084540

# Row 12
$ws.Range("D12").Value = -1205334
$ws.Range("E12").Value = -1974746
$ws.Range("F12").Value = -3237613
$ws.Range("G12").Value = -4409368
$ws.Range("H12").Value = -6968838

# Row 13
$ws.Range("D13").Value = 884968
$ws.Range("E13").Value = 1580093
$ws.Range("F13").Value = 4052848
$ws.Range("G13").Value = 5157558
$ws.Range("H13").Value = 6115702

# Row 14
$ws.Range("D14").Value = -162955
$ws.Range("E14").Value = -228431
$ws.Range("F14").Value = -334172
$ws.Range("G14").Value = -456236
$ws.Range("H14").Value = -596407

# Row 15
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0

# Row 16
$ws.Range("D16").Value = 42595
$ws.Range("E16").Value = 135151
$ws.Range("F16").Value = 134884
$ws.Range("G16").Value = 32241
$ws.Range("H16").Value = 328860

# Row 17
$ws.Range("D17").Value = 764608
$ws.Range("E17").Value = 1486813
$ws.Range("F17").Value = 3853560
$ws.Range("G17").Value = 4733563
$ws.Range("H17").Value = 5848155

# Row 18
$ws.Range("D18").Value = -174779
$ws.Range("E18").Value = -179671
$ws.Range("F18").Value = -70702
$ws.Range("G18").Value = -146406
$ws.Range("H18").Value = -461471

# Row 19
$ws.Range("D19").Value = 8038
$ws.Range("E19").Value = -66581
$ws.Range("F19").Value = 51876
$ws.Range("G19").Value = 81079
$ws.Range("H19").Value = 48199

# Row 20
$ws.Range("D20").Value = 597867
$ws.Range("E20").Value = 1240561
$ws.Range("F20").Value = 3834734
$ws.Range("G20").Value = 4668236
$ws.Range("H20").Value = 5434883

# Row 21
$ws.Range("D21").Value = -47032
$ws.Range("E21").Value = -120469
$ws.Range("F21").Value = -241427
$ws.Range("G21").Value = -576085
$ws.Range("H21").Value = -757110

# Row 22
$ws.Range("D22").Value = 550835
$ws.Range("E22").Value = 1120092
$ws.Range("F22").Value = 3593307
$ws.Range("G22").Value = 4092151
$ws.Range("H22").Value = 4677773

# Row 23
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0

# Row 24
$ws.Range("D24").Value = 550835
$ws.Range("E24").Value = 1120092
$ws.Range("F24").Value = 3593307
$ws.Range("G24").Value = 4092151
$ws.Range("H24").Value = 4677773

# Row 25
$ws.Range("D25").Value = 1102
$ws.Range("E25").Value = 2240
$ws.Range("F25").Value = 5026
$ws.Range("G25").Value = 5723
$ws.Range("H25").Value = 6542

# Row 26
$ws.Range("D26").Value = 500000
$ws.Range("E26").Value = 500000
$ws.Range("F26").Value = 715000
$ws.Range("G26").Value = 715000
$ws.Range("H26").Value = 715000

# Row 27
$ws.Range("D27").Value = 770
$ws.Range("E27").Value = 1567
$ws.Range("F27").Value = 5026
$ws.Range("G27").Value = 5723
$ws.Range("H27").Value = 6542

